# edit.ps1 — reproduces the two semantic changes from the target diff:
#
# 1. The table on slide 6 (the "SOURCES OF FINANCE" table) gets its table
#    style switched from the custom "{06FDC9FD-1CC6-4EDE-AB5A-B6F676348F7A}"
#    style to the built-in "{44383B47-7BF8-4F54-8286-7F6D5ED7B461}" style.
#
# 2. The presentation's theme (ppt/theme/theme1.xml, used by the Slide
#    Master) changes from the custom "Integral" palette to the stock
#    "Office" palette (the deck's Notes Master already carried that
#    "Office Theme" palette as its own theme, so after this edit the two
#    theme parts in the package end up holding each other's former colors).

$p = $ppt.ActivePresentation

# --- 1. Update the table style on slide 6 ---------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{44383B47-7BF8-4F54-8286-7F6D5ED7B461}")
    }
}

# --- 2. Re-colour the presentation theme to the stock "Office" palette ----
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# Order matches the standard 12-slot DrawingML colour scheme:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeColors = @(
    0,           # dk1      000000
    16777215,    # lt1      FFFFFF
    6968388,     # dk2      44546A
    15132391,    # lt2      E7E6E6
    13998939,    # accent1  5B9BD5
    3243501,     # accent2  ED7D31
    10855845,    # accent3  A5A5A5
    49407,       # accent4  FFC000
    12874308,    # accent5  4472C4
    4697456,     # accent6  70AD47
    12673797,    # hlink    0563C1
    7491477      # folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i - 1]
}
